$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 361, shifting rows 361-425 down to 362-426
$ws.Rows.Item(361).Insert()

# Populate the new row 361 with the new pineapple price record
$ws.Cells.Item(361, 1).Value = 10
$ws.Cells.Item(361, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(361, 3).Value = "La Araucanía"
$ws.Cells.Item(361, 4).Value = 44637
$ws.Cells.Item(361, 5).Value = 9
$ws.Cells.Item(361, 6).Value = "Fruta"
$ws.Cells.Item(361, 7).Value = 100108
$ws.Cells.Item(361, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(361, 9).Value = 100108005
$ws.Cells.Item(361, 10).Value = "Piña"
$ws.Cells.Item(361, 11).Value = "Caramelo"
$ws.Cells.Item(361, 12).Value = "Primera"
$ws.Cells.Item(361, 13).Value = 155
$ws.Cells.Item(361, 14).Value = 18000
$ws.Cells.Item(361, 15).Value = 19000
$ws.Cells.Item(361, 16).Value = 18581
$ws.Cells.Item(361, 17).Value = "$/caja 12 unidades"
$ws.Cells.Item(361, 18).Value = "Ecuador"
$ws.Cells.Item(361, 19).Value = 1548
$ws.Cells.Item(361, 20).Value = 12
